$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 5252
$ws.Range("E3").Value = 6395
$ws.Range("E4").Value = 8107
$ws.Range("E5").Value = 3718
$ws.Range("E6").Value = 6009
$ws.Range("E7").Value = 13877
$ws.Range("E8").Value = 15603
$ws.Range("E9").Value = 11396
$ws.Range("E10").Value = 19883
$ws.Range("E11").Value = 6974
$ws.Range("E12").Value = 19523
$ws.Range("E13").Value = 8218
